$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("D5").Value = "2016-02-29 12:55:08"
$wsDeDe.Range("D5").Value = "2016-02-29 12:55:18"
